$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated SFH UFA baseline model output values (column B, rows 2-452)
# Values recalibrated per commit: "calibration of energy use modeling by renovation level"
$values = @{
    2 = 24182101.97820457
    3 = 29103.36738456176
    4 = 29545.63453354353
    5 = 29997.31595861042
    6 = 30458.73678669694
    7 = 30930.23515986492
    8 = 31412.16259442532
    9 = 31904.88433427436
    10 = 32408.77969701503
    11 = 32924.24241093556
    12 = 33451.68094157633
    13 = 33991.51880618763
    14 = 34544.19487386242
    15 = 35110.16365008971
    16 = 35689.89554318571
    17 = 36283.87711115431
    18 = 36892.61128641399
    19 = 37516.61757657271
    20 = 38156.43223899048
    21 = 38812.60842664915
    22 = 39485.71630347187
    23 = 40176.34312625887
    24 = 40885.09329116678
    25 = 41612.58834233681
    26 = 42359.46693993507
    27 = 43126.38478558441
    28 = 43914.01450238063
    29 = 44723.04546722399
    30 = 45554.18359312595
    31 = 46408.15105881709
    32 = 47285.68598358224
    33 = 48187.54204489486
    34 = 49114.4880365626
    35 = 50067.30736535523
    36 = 51046.79748387059
    37 = 52053.76925787541
    38 = 53089.04626596631
    39 = 54153.46403020087
    40 = 55247.86917584204
    41 = 56373.11851892914
    42 = 57530.07808026388
    43 = 58719.62202523054
    44 = 59942.6315279113
    45 = 61199.99355956532
    46 = 62492.59960063724
    47 = 63821.3442762475
    48 = 65187.12391564368
    49 = 66590.83503544162
    50 = 68033.37274779624
    51 = 69515.62909447624
    52 = 71038.49130789773
    53 = 72602.84000108093
    54 = 74209.54728820956
    55 = 75859.47483854993
    56 = 77553.47186584784
    57 = 79292.3730566809
    58 = 81076.99644117728
    59 = 82908.14120942533
    60 = 84786.58547821651
    61 = 86713.08401246098
    62 = 88688.36590605169
    63 = 90713.13222746048
    64 = 92788.05363602587
    65 = 94913.76797414836
    66 = 97090.87784272023
    67 = 99319.9481654433
    68 = 101601.5037497326
    69 = 103936.0268512889
    70 = 106323.9547498022
    71 = 108765.677343731
    72 = 111261.5347725448
    73 = 113811.8150743027
    74 = 116416.7518876409
    75 = 119076.5222066312
    76 = 121791.2441975429
    77 = 124560.975086528
    78 = 127385.7091274912
    79 = 130265.3756591299
    80 = 133199.8372606734
    81 = 136188.8880153309
    82 = 139232.2518907566
    83 = 142329.5812456139
    84 = 145480.4554712684
    85 = 148684.3797773647
    86 = 151940.7841300619
    87 = 155249.0223509393
    88 = 158608.3713852452
    89 = 162018.0307469041
    90 = 165477.1221473856
    91 = 168984.6893161852
    92 = 172539.6980188388
    93 = 176141.0362785915
    94 = 179787.5148075092
    95 = 183477.8676517937
    96 = 187210.7530554775
    97 = 190984.7545465426
    98 = 194798.3822482077
    99 = 198650.0744178414
    100 = 202538.1992150421
    101 = 206461.0566995505
    102 = 210416.881059512
    103 = 214403.843068734
    104 = 218420.0527716048
    105 = 222463.5623936387
    106 = 226532.3694735644
    107 = 230624.4202137394
    108 = 234737.6130431476
    109 = 238869.8023877447
    110 = 243018.8026411888
    111 = 247182.3923283959
    112 = 251358.3184538406
    113 = 255544.3010250417
    114 = 259738.037741738
    115 = 263937.2088396235
    116 = 268139.4820773657
    117 = 272342.5178546925
    118 = 276543.9744487477
    119 = 280741.5133549829
    120 = 284932.8047187958
    121 = 289115.532843054
    122 = 293287.4017564809
    123 = 297446.1408272049
    124 = 301589.5104055835
    125 = 305715.3074799349
    126 = 309821.3713286376
    127 = 313905.5891516459
    128 = 317965.9016647415
    129 = 322000.3086390879
    130 = 326006.8743695837
    131 = 329983.7330542579
    132 = 333929.0940686666
    133 = 337841.2471180224
    134 = 341718.5672505762
    135 = 345559.519716844
    136 = 349362.664657784
    137 = 353126.6616078786
    138 = 356850.2737975773
    139 = 360532.3722410263
    140 = 364171.9395960031
    141 = 367768.0737827354
    142 = 371319.9913496976
    143 = 374827.0305754459
    144 = 378288.6542954934
    145 = 381704.4524454235
    146 = 385074.1443110013
    147 = 388397.5804785439
    148 = 391674.7444780828
    149 = 394905.7541146747
    150 = 398090.8624830688
    151 = 401230.4586621471
    152 = 404325.068087659
    153 = 407375.3526013208
    154 = 410382.1101771028
    155 = 413346.2743252177
    156 = 416268.9131763856
    157 = 419151.2282501101
    158 = 421994.5529105499
    159 = 424800.3505166919
    160 = 427570.2122725116
    161 = 430305.8547853751
    162 = 433009.1173410746
    163 = 435681.9589051903
    164 = 438326.4548614437
    165 = 440944.7934980365
    166 = 443539.2722547883
    167 = 446112.2937437967
    168 = 448666.3615571994
    169 = 451204.0758766997
    170 = 453728.1288998845
    171 = 456241.3000986105
    172 = 458746.4513257309
    173 = 461246.5217865184
    174 = 463744.5228916048
    175 = 466243.5330083659
    176 = 468746.6921282437
    177 = 471257.1964673362
    178 = 473778.2930176578
    179 = 476313.2740667704
    180 = 478865.4717030184
    181 = 481438.2523239235
    182 = 484035.0111647897
    183 = 486659.1668643216
    184 = 489314.1560838788
    185 = 492003.4281968626
    186 = 494730.4400635625
    187 = 497498.6509070721
    188 = 500311.5173052275
    189 = 503172.4883121713
    190 = 506085.0007240582
    191 = 509052.4745010047
    192 = 512078.3083578672
    193 = 515165.8755358686
    194 = 518318.5197648237
    195 = 521539.5514267895
    196 = 524832.2439301702
    197 = 528199.8303023838
    198 = 531645.5000088216
    199 = 535172.3960049057
    200 = 538783.6120267219
    201 = 542482.1901256988
    202 = 546271.1184511291
    203 = 550153.4520478038
    204 = 554131.8335083345
    205 = 558209.1888494854
    206 = 562388.2735527877
    207 = 566671.7819179274
    208 = 571062.3463729275
    209 = 575562.5370383358
    210 = 580174.8615429065
    211 = 584901.7650879147
    212 = 589745.6307565001
    213 = 594708.7800635347
    214 = 599793.4737413282
    215 = 605001.9127559266
    216 = 610336.2395473893
    217 = 615798.539488025
    218 = 621390.8425515642
    219 = 627115.1251854961
    220 = 632973.3123790545
    221 = 638967.2799191637
    222 = 645098.8568240642
    223 = 651369.8279484758
    224 = 657781.9367486307
    225 = 664336.8882002792
    226 = 671036.3518575506
    227 = 677881.9650464065
    228 = 684875.3361796064
    229 = 692018.0481857682
    230 = 699311.6620414359
    231 = 706757.7203972075
    232 = 714357.7512865379
    233 = 722113.2719088882
    234 = 730025.7924767904
    235 = 738096.8201162544
    236 = 746327.8628121845
    237 = 754720.4333880187
    238 = 763276.0535104405
    239 = 771996.2577099768
    240 = 780882.5974078375
    241 = 789936.6449394722
    242 = 799159.9975669304
    243 = 808554.2814698762
    244 = 818121.1557067464
    245 = 827862.3161382569
    246 = 837779.4993032403
    247 = 847874.4862399215
    248 = 858149.1062431417
    249 = 868605.2405499992
    250 = 879244.8259455865
    251 = 890069.8582811452
    252 = 901082.3958958519
    253 = 912284.562935676
    254 = 923678.5525605899
    255 = 935266.6300320991
    256 = 947051.1356744417
    257 = 959034.4877008732
    258 = 971219.1848971776
    259 = 983607.8091549231
    260 = 996203.027847002
    261 = 1009007.596036025
    262 = 1022024.358509763
    263 = 1035256.251633283
    264 = 1048706.305010542
    265 = 1062377.642948413
    266 = 1076273.485710652
    267 = 1090397.150558132
    268 = 1104752.05256271
    269 = 1119341.705186585
    270 = 1134169.720619095
    271 = 1149239.809859537
    272 = 1164555.782538138
    273 = 1180121.546463526
    274 = 1195941.106887897
    275 = 1212018.565478925
    276 = 1228358.118987373
    277 = 1244964.057600851
    278 = 1261840.762971281
    279 = 1278992.705905907
    280 = 1296424.443709748
    281 = 1314140.617168707
    282 = 1332145.947160923
    283 = 1350445.230885123
    284 = 1369043.337693794
    285 = 1387945.204519483
    286 = 1407155.830882689
    287 = 1426680.273468147
    288 = 1446523.640259787
    289 = 1466691.084221005
    290 = 1487187.796509924
    291 = 1508018.999218091
    292 = 1529189.937621736
    293 = 1550705.871935356
    294 = 1572572.068558673
    295 = 1594793.790805205
    296 = 1617376.289106994
    297 = 1640324.79068555
    298 = 1663644.488682524
    299 = 1687340.530744179
    300 = 1711418.007055671
    301 = 1735881.937820126
    302 = 1760737.260180198
    303 = 1233523.122231114
    304 = 1251240.462552959
    305 = 1269237.913509074
    306 = 1287518.557543174
    307 = 1306085.368823881
    308 = 1324941.201935635
    309 = 1344088.780209123
    310 = 1363530.683699444
    311 = 1383269.336822803
    312 = 1403306.995667197
    313 = 1423645.734990872
    314 = 1444287.434925416
    315 = 1465233.76740532
    316 = 1486486.182344939
    317 = 1508045.893587124
    318 = 1529913.864652751
    319 = 1552090.794318502
    320 = 1574577.10205702
    321 = 1597372.91337597
    322 = 1620478.045091042
    323 = 1643891.990577831
    324 = 1667613.905042759
    325 = 1691642.5908644
    326 = 1715976.48305057
    327 = 1740613.634866747
    328 = 1765551.703690816
    329 = 1790787.937151277
    330 = 1816319.159610766
    331 = 1842141.759058074
    332 = 1868251.67447446
    333 = 1894644.383741266
    334 = 1921314.892160501
    335 = 1948257.721658577
    336 = 1975466.900747313
    337 = 2002935.95531716
    338 = 2030657.900337348
    339 = 2058625.232539059
    340 = 2086829.924160938
    341 = 2115263.417829871
    342 = 2143916.622655388
    343 = 2172779.911613613
    344 = 2201843.120293225
    345 = 2231095.547078443
    346 = 2260525.954835607
    347 = 2290122.574177071
    348 = 2319873.108360208
    349 = 2349764.739889903
    350 = 2379784.138874648
    351 = 2409917.473195491
    352 = 2440150.420529897
    353 = 2470468.182276341
    354 = 2500855.499412627
    355 = 2531296.6703198
    356 = 2561775.570592089
    357 = 2592275.674848779
    358 = 2504897.739758637
    359 = 2534018.737362935
    360 = 2563110.577374399
    361 = 2592156.191963608
    362 = 2621138.263137701
    363 = 2650039.254896829
    364 = 2678841.447177625
    365 = 2707526.971524006
    366 = 2736077.848424881
    367 = 2764476.026248843
    368 = 2792703.421688373
    369 = 2820741.961630428
    370 = 2848573.626347312
    371 = 2876180.493904503
    372 = 2903544.785668758
    373 = 3125235.551509051
    374 = 3153843.371861015
    375 = 3182137.049779046
    376 = 3210098.741399611
    377 = 3237711.011024292
    378 = 3264956.883658133
    379 = 3291819.897653168
    380 = 3318284.157283706
    381 = 3344334.38507777
    382 = 3369955.973725168
    383 = 3667323.400436685
    384 = 3694028.905048542
    385 = 3720228.96444211
    386 = 3745911.240577997
    387 = 3771064.386028036
    388 = 3795678.093251877
    389 = 3819743.142082406
    390 = 3843251.445242418
    391 = 3866196.091717481
    392 = 3888571.387820199
    393 = 3300971.426083894
    394 = 3318888.313798556
    395 = 3336316.641165851
    396 = 3353256.015285067
    397 = 3369707.210392496
    398 = 3385672.188180401
    399 = 3401154.115344614
    400 = 3416157.378261506
    401 = 3430687.59471599
    402 = 3444751.622589133
    403 = 2380099.904857051
    404 = 2389154.917928257
    405 = 2397908.399420185
    406 = 2406368.499116863
    407 = 2414544.21709316
    408 = 2422445.396156697
    409 = 2430082.711891975
    410 = 2437467.660259086
    411 = 2444612.5427203
    412 = 2451530.448846592
    413 = 2555346.58376493
    414 = 2562109.883078001
    415 = 2568682.754103084
    416 = 2575081.846014949
    417 = 2581324.496218503
    418 = 2587428.69572821
    419 = 2593413.051640937
    420 = 2599296.746655782
    421 = 2605099.495597587
    422 = 2610841.498896472
    423 = 2616543.392986387
    424 = 2622226.197574971
    425 = 2627911.259756727
    426 = 2633620.19492895
    427 = 2639374.824492062
    428 = 2645197.110315716
    429 = 2651109.085961381
    430 = 2657132.784671591
    431 = 2663290.164143504
    432 = 2669603.028123542
    433 = 2676092.944883685
    434 = 2682781.162653679
    435 = 2689688.522119531
    436 = 2696835.366112857
    437 = 2704241.446662311
    438 = 2711925.829594668
    439 = 2719906.7969263
    440 = 2728201.747307728
    441 = 2736827.094840323
    442 = 2745798.16661708
    443 = 2755129.09938701
    444 = 2764832.735793971
    445 = 2774920.520682252
    446 = 2785402.398013207
    447 = 2796286.708984035
    448 = 2807580.09198821
    449 = 2819287.385107951
    450 = 2831411.531859202
    451 = 2843953.490973537
    452 = 2856912.151019436
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
